$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Shrub" -> "CSS" (Ashish's "Shrub" becomes "CSS" to match enzyme activity data)
# A2 was "Grassland" -> becomes "CSS"
# B2 was "Shrub" -> becomes "Grassland"
$ws.Range("A2").Value = "CSS"
$ws.Range("B2").Value = "Grassland"

# Update statistic values: meandiff, lower, upper flip sign
$ws.Range("C2").Value = -1.2627
$ws.Range("E2").Value = -1.5765
$ws.Range("F2").Value = -0.9489
